# Applies the numeric updates from the scheduled-runner data refresh
# (currentAveragePrice* / LevePrice* / LeveProfit* columns H:N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit sheets.
$wb = $excel.ActiveWorkbook

# ALC row 39
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 2029.25
$ws.Range("I39").Value = 283.2
$ws.Range("J39").Value = 3276.4285
$ws.Range("K39").Value = 849.5999999999999
$ws.Range("L39").Value = 9829.2855
$ws.Range("M39").Value = -553.5999999999999
$ws.Range("N39").Value = -10421.2855

# ALC row 64
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 10668
$ws.Range("J64").Value = 11777.5
$ws.Range("L64").Value = 11777.5
$ws.Range("N64").Value = -12273.5

# ALC row 67
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 10668
$ws.Range("J67").Value = 11777.5
$ws.Range("L67").Value = 11777.5
$ws.Range("N67").Value = -13493.5

# ALC row 87
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 19000
$ws.Range("I87").Value = 8500
$ws.Range("K87").Value = 8500
$ws.Range("M87").Value = -7252

# ALC row 90
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H90").Value = 19000
$ws.Range("I90").Value = 8500
$ws.Range("K90").Value = 25500
$ws.Range("M90").Value = -19260

# ALC row 125
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 2422.9167
$ws.Range("I125").Value = 1110.125
$ws.Range("K125").Value = 9991.125
$ws.Range("M125").Value = -7531.125

# ALC row 131
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 7779
$ws.Range("I131").Value = 3031.6667
$ws.Range("K131").Value = 9095.000100000001
$ws.Range("M131").Value = -4055.000100000001

# ALC row 135
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 14085.485
$ws.Range("I135").Value = 15274.678
$ws.Range("J135").Value = 4869.25
$ws.Range("K135").Value = 137472.102
$ws.Range("L135").Value = 43823.25
$ws.Range("M135").Value = -134937.102
$ws.Range("N135").Value = -48893.25

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3673.7017
$ws.Range("I138").Value = 5111.1763
$ws.Range("J138").Value = 3062.775
$ws.Range("K138").Value = 15333.5289
$ws.Range("L138").Value = 9188.325000000001
$ws.Range("M138").Value = -10193.5289
$ws.Range("N138").Value = -19468.325

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4305.875
$ws.Range("I61").Value = 2730
$ws.Range("K61").Value = 2730
$ws.Range("M61").Value = -2518

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 4305.875
$ws.Range("I136").Value = 2730
$ws.Range("K136").Value = 8190
$ws.Range("M136").Value = -5640

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4904.387
$ws.Range("I86").Value = 4843.5
$ws.Range("J86").Value = 5015.091
$ws.Range("K86").Value = 4843.5
$ws.Range("L86").Value = 5015.091
$ws.Range("M86").Value = -3720.5
$ws.Range("N86").Value = -7261.091

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 4904.387
$ws.Range("I89").Value = 4843.5
$ws.Range("J89").Value = 5015.091
$ws.Range("K89").Value = 24217.5
$ws.Range("L89").Value = 25075.455
$ws.Range("M89").Value = -18601.5
$ws.Range("N89").Value = -36307.455

# BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3075.652
$ws.Range("I105").Value = 2693.3076
$ws.Range("J105").Value = 3572.7
$ws.Range("K105").Value = 2693.3076
$ws.Range("L105").Value = 3572.7
$ws.Range("M105").Value = -946.3076000000001
$ws.Range("N105").Value = -7066.7

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2219.6177
$ws.Range("I31").Value = 2306.077
$ws.Range("J31").Value = 2166.0952
$ws.Range("K31").Value = 2306.077
$ws.Range("L31").Value = 2166.0952
$ws.Range("M31").Value = -2011.077
$ws.Range("N31").Value = -2756.0952

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2219.6177
$ws.Range("I34").Value = 2306.077
$ws.Range("J34").Value = 2166.0952
$ws.Range("K34").Value = 2306.077
$ws.Range("L34").Value = 2166.0952
$ws.Range("M34").Value = -2104.077
$ws.Range("N34").Value = -2570.0952

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2835.625
$ws.Range("I132").Value = 2527.8
$ws.Range("K132").Value = 7583.400000000001
$ws.Range("M132").Value = -5053.400000000001

# CUL row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 897
$ws.Range("I5").Value = 697
$ws.Range("K5").Value = 2091
$ws.Range("M5").Value = -1979

# CUL row 22
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 23199.8
$ws.Range("J22").Value = 28499.75
$ws.Range("L22").Value = 85499.25
$ws.Range("N22").Value = -85837.25

# CUL row 27
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H27").Value = 23199.8
$ws.Range("J27").Value = 28499.75
$ws.Range("L27").Value = 85499.25
$ws.Range("N27").Value = -85703.25

# CUL row 35
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()

# CUL row 39
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 10166.333

# CUL row 40
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 81
$ws.Range("I40").Value = 74.833336
$ws.Range("J40").Value = 93.333336
$ws.Range("K40").Value = 299.333344
$ws.Range("L40").Value = 373.333344
$ws.Range("M40").Value = -230.333344
$ws.Range("N40").Value = -511.333344

# CUL row 46
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 1620.875
$ws.Range("I46").Value = 774.5
$ws.Range("J46").Value = 1790.15
$ws.Range("K46").Value = 2323.5
$ws.Range("L46").Value = 5370.450000000001
$ws.Range("M46").Value = -2232.5
$ws.Range("N46").Value = -5552.450000000001

# CUL row 54
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H54").Value = 10333
$ws.Range("J54").Value = 10333
$ws.Range("L54").Value = 30999
$ws.Range("N54").Value = -32117

# CUL row 58
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 9374.75

# CUL row 59
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value = 4224.75
$ws.Range("I59").Value = 1950
$ws.Range("K59").Value = 5850
$ws.Range("M59").Value = -5310

# CUL row 76
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H76").Value = 12428.286
$ws.Range("J76").Value = 14000
$ws.Range("L76").Value = 42000
$ws.Range("N76").Value = -42766

# CUL row 79
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H79").Value = 12428.286
$ws.Range("J79").Value = 14000
$ws.Range("L79").Value = 42000
$ws.Range("N79").Value = -44652

# CUL row 129
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 8335674.5
$ws.Range("I129").Value = 489.14285
$ws.Range("J129").Value = 12823851
$ws.Range("K129").Value = 1467.42855
$ws.Range("L129").Value = 38471553
$ws.Range("M129").Value = 3532.57145
$ws.Range("N129").Value = -38481553

# CUL row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 897
$ws.Range("I135").Value = 697
$ws.Range("K135").Value = 6273
$ws.Range("M135").Value = -3738

# CUL row 140
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 6587075
$ws.Range("I140").Value = 11906184
$ws.Range("J140").Value = 16411.059
$ws.Range("K140").Value = 35718552
$ws.Range("L140").Value = 49233.177
$ws.Range("M140").Value = -35713372
$ws.Range("N140").Value = -59593.177

# GSM row 107
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 746.2727
$ws.Range("I107").Value = 792.1
$ws.Range("J107").Value = 288
$ws.Range("K107").Value = 792.1
$ws.Range("L107").Value = 288
$ws.Range("M107").Value = 1127.9
$ws.Range("N107").Value = -4128

# LTW row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2491.75
$ws.Range("J68").Value = 2475
$ws.Range("L68").Value = 2475
$ws.Range("N68").Value = -3973

# LTW row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 2491.75
$ws.Range("J71").Value = 2475
$ws.Range("L71").Value = 12375
$ws.Range("N71").Value = -19863

# WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 4714.4287
$ws.Range("I126").Value = 3500.2727
$ws.Range("K126").Value = 10500.8181
$ws.Range("M126").Value = -8030.8181

# WVR row 139
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H139").Value = 68514.81
$ws.Range("J139").Value = 68514.81
$ws.Range("L139").Value = 68514.81
$ws.Range("N139").Value = -78794.81
